# Apply angle-bracket wrapping to @base/@prefix URI values on the "Sheet4" metadata sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet4")

$ws.Range("D7").Value = "<http://junk/just-for-fun>"
$ws.Range("D1").Value = "<http://seman.tc/data/northwind>"
$ws.Range("D2").Value = "<http://seman.tc/models/northwind#>"
$ws.Range("D3").Value = "<http://schema.org/>"
$ws.Range("D4").Value = "<http://xmlns.com/foaf/0.1/>"
$ws.Range("D5").Value = "<http://seman.tc/models/northwind#>"
$ws.Range("D6").Value = "<http://purl.org/dc/terms/>"

# Widen column D to fit the new, longer bracketed values (best-fit sized for the content).
$ws.Columns.Item(4).ColumnWidth = 35.7

# Move the active selection to D7.
$ws.Range("D7").Select()
